# Refresh the per-minute crypto snapshot (price + 1h volume change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.825.49'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.633.14'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '214.74'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").Value = '19.91'
$ws.Range("E10").Value = '  +2.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '1.659.92'
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").Value = '1.858.52'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '63.04'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '25.830.03'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '194.02'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("E21").Value = '  +2.16%  '
$ws.Range("D22").Value = '9.92'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").Value = '6.19'
$ws.Range("E23").Value = '  +3.14%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -2.88%  '
$ws.Range("D26").Value = '139.74'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("E27").Value = '  -2.67%  '
$ws.Range("D28").Value = '6.83'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '0.0494'
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("E33").Value = '  +2.36%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '0.551'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").Value = '1.121.39'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("D43").Value = '99.72'
$ws.Range("E43").Value = '  +2.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.800'
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").Value = '0.0₆0108'
$ws.Range("E45").Value = '  -3.63%  '
$ws.Range("D46").Value = '55.45'
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").Value = '  +7.75%  '
$ws.Range("E51").Value = '  -0.24%  '
